$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# The "meta" sheet stores key/value pairs, one per row, with a trailing
# blank (but styled) row acting as a spacer at the bottom. We are adding a
# new "style" = "default" key/value pair. Insert a fresh row above the
# existing blank spacer row (pushing the spacer down one row) and fill the
# newly inserted row with the new key/value pair.
$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = "style"
$ws.Range("B7").Value = "default"
